$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.756.48'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '3.057.36'
$ws.Range('E3').Value = '  -1.75%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'536.50"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.89%  '
$ws.Range('D6').Value = "'132.43"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.80%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.048.98'
$ws.Range('E8').Value = '  -1.78%  '
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('D10').Value = "'0.154"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.15%  '
$ws.Range('E11').Value = '  -10.24%  '
$ws.Range('D12').Value = "'0.450"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.04%  '
$ws.Range('D13').Value = "'0.0000223"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('D14').Value = "'34.07"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.24%  '
$ws.Range('D15').Value = '3.552.86'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').Value = '62.787.97'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = '3.059.34'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').Value = "'6.61"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').Value = "'480.58"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.40%  '
$ws.Range('E21').Value = '  -3.70%  '
$ws.Range('D22').Value = "'0.691"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.37%  '
$ws.Range('D23').Value = "'7.08"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.21%  '
$ws.Range('D24').Value = "'78.87"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').Value = "'12.03"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.91%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = "'2.70"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.47%  '
$ws.Range('D28').Value = "'8.03"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.80%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = "'25.86"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.08%  '
$ws.Range('E31').Value = '  -9.55%  '
$ws.Range('E32').Value = '  -2.12%  '
$ws.Range('B33').Value = 'Stacks'
$ws.Range('C33').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D33').Value = "'2.35"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.75%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = "'56.82"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.34%  '
$ws.Range('D35').Value = "'5.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.14%  '
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').Value = "'481.57"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -10.58%  '
$ws.Range('D38').Value = '3.094.39'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').Value = "'0.0393"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.11%  '
$ws.Range('D40').Value = "'0.0792"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('E41').Value = '  -3.78%  '
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('D43').Value = "'2.60"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.58%  '
$ws.Range('E44').Value = '  -3.51%  '
$ws.Range('E46').Value = '  +5.32%  '
$ws.Range('D47').Value = "'121.22"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('E48').Value = '  -6.15%  '
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = "'2.00"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.28%  '
